$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Automation"
$ws.Range("A4").Value = "php"
$ws.Range("A5").Value = "JS"
$ws.Range("A6").Value = "Ajax"

$ws.Range("A7").Select()
